# Release plan and notes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Release Plan")

# --- Insert a new column D "Actual Release Date" by copying column C's formatting ---
$ws.Columns.Item(3).Copy()
$ws.Columns.Item(4).Insert()

# Header
$ws.Range("C1").Value = "Planned Release Date"
$ws.Range("D1").Value = "Actual Release Date"

# Fix typo in shared text "Edit Profile- Srudent/Staff" -> "Edit Profile- Student/Staff"
$ws.Range("B8").Value = "Edit Profile- Student/Staff"

# New column D values (Actual Release Date) - merged similarly to column C
$ws.Range("D2").Value = 42463
$ws.Range("D4").Value = 42467
$ws.Range("D6").Value = 42471
$ws.Range("D9").Value = 42475

# Merge D column cells to match C column merges
$ws.Range("D2:D3").Merge()
$ws.Range("D4:D5").Merge()
$ws.Range("D6:D8").Merge()
$ws.Range("D9:D10").Merge()

# Column D width
$ws.Columns.Item(4).ColumnWidth = 17.33203125

# Row 1 height (header row wraps to 2 lines now)
$ws.Rows.Item(1).RowHeight = 30

# Update selection to match target
$ws.Range("B8").Select()

# Window view size changes
$excel.ActiveWindow.WindowState = -4143
